# Apply "5-year growth is frozen before first drop in projections" logic to the
# GroupGrowthRateMax sheet.
#
# Each group of 6 consecutive data rows (periods 2025, 2030, 2035, 2040, 2045, 2050)
# shares the same region/group_name. From the period-2030 row onward, once growth
# starts decreasing in the original projection, the rate is "frozen" at the 2030
# value for all subsequent periods (2035, 2040, 2045, 2050).
#
# In addition, every note in column E gets an extra clause inserted right after
# "...logistic diffusion model" (before the trailing "; Median growth scenario...").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GroupGrowthRateMax")

$oldPhrase = "logistic diffusion model;"
$newPhrase = "logistic diffusion model, 5-year growth is frozen before first drop in projections;"

$firstDataRow = 2
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt $firstDataRow) {
    $lastRow = 145
}

for ($start = $firstDataRow; $start -le $lastRow; $start += 6) {
    $frozenValue = $ws.Cells.Item($start + 1, 4).Value2

    for ($offset = 2; $offset -le 5; $offset++) {
        $row = $start + $offset
        if ($row -le $lastRow) {
            $ws.Cells.Item($row, 4).Value2 = $frozenValue
        }
    }
}

for ($row = $firstDataRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $note = $cell.Value2
    if ($note -ne $null -and $note.ToString().Contains($oldPhrase)) {
        $cell.Value2 = $note.ToString().Replace($oldPhrase, $newPhrase)
    }
}
